# Weekly update: shift existing history down by 3 rows (insert 3 new rows
# right before row 423) and populate the new top rows with this week's
# "Calameño" entries (Extra / Primera / Segunda) for the Vega Modelo de
# Temuco - Melón sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 423; everything from old row 423 downward
# (through old row 522) shifts down to rows 426-525, and the sheet
# dimension grows from A1:R522 to A1:R525.
$ws.Range("A423:R425").Insert()

# New row 423: Melón / Calameño / Extra
$ws.Cells.Item(423, 4).Value2 = 44588
$ws.Cells.Item(423, 9).Value2 = "Extra"
$ws.Cells.Item(423, 10).Value2 = 1000
$ws.Cells.Item(423, 11).Value2 = 1200
$ws.Cells.Item(423, 12).Value2 = 1200
$ws.Cells.Item(423, 13).Value2 = 1200
$ws.Cells.Item(423, 15).Value2 = "Región del Maule"
$ws.Cells.Item(423, 16).Value2 = 1200

# New row 424: Melón / Calameño / Primera
$ws.Cells.Item(424, 4).Value2 = 44588
$ws.Cells.Item(424, 9).Value2 = "Primera"
$ws.Cells.Item(424, 10).Value2 = 10000
$ws.Cells.Item(424, 11).Value2 = 1000
$ws.Cells.Item(424, 12).Value2 = 1000
$ws.Cells.Item(424, 13).Value2 = 1000
$ws.Cells.Item(424, 15).Value2 = "Región del Maule"
$ws.Cells.Item(424, 16).Value2 = 1000

# New row 425: Melón / Calameño / Segunda
$ws.Cells.Item(425, 4).Value2 = 44588
$ws.Cells.Item(425, 9).Value2 = "Segunda"
$ws.Cells.Item(425, 10).Value2 = 2000
$ws.Cells.Item(425, 11).Value2 = 800
$ws.Cells.Item(425, 12).Value2 = 800
$ws.Cells.Item(425, 13).Value2 = 800
$ws.Cells.Item(425, 15).Value2 = "Región del Maule"
$ws.Cells.Item(425, 16).Value2 = 800

# Shared fields that are the same on every row of this block and need to be
# filled in for the 3 freshly inserted rows (A/B/C/E/F/G/H/N/Q/R).
foreach ($r in 423..425) {
    $ws.Cells.Item($r, 1).Value2 = 10
    $ws.Cells.Item($r, 2).Value2 = "Vega Modelo de Temuco"
    $ws.Cells.Item($r, 3).Value2 = "La Araucanía"
    $ws.Cells.Item($r, 5).Value2 = 9
    $ws.Cells.Item($r, 6).Value2 = 100112027
    $ws.Cells.Item($r, 7).Value2 = "Melón"
    $ws.Cells.Item($r, 8).Value2 = "Calameño"
    $ws.Cells.Item($r, 14).Value2 = "$/unidad"
    $ws.Cells.Item($r, 17).Value2 = 1
    $ws.Cells.Item($r, 18).Value2 = "Hortaliza"
}
